# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# 104a7d0d-... entry (row 3) on both the "zh-cn" and "de-de" sheets.
# Rows 3 and 4 share the same underlying values, so both rows must be
# updated to keep them in sync.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E3").Value = "2016-03-20 04:17:32"
$ws_zhcn.Range("E4").Value = "2016-03-20 04:17:32"
$ws_zhcn.Range("H3").Value = "2016-03-20 04:17:52"
$ws_zhcn.Range("H4").Value = "2016-03-20 04:17:52"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E3").Value = "2016-03-20 04:17:36"
$ws_dede.Range("E4").Value = "2016-03-20 04:17:36"
$ws_dede.Range("H3").Value = "2016-03-20 04:17:58"
$ws_dede.Range("H4").Value = "2016-03-20 04:17:58"
